$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.960.28"
$ws.Range("E2").Value = "  -0.50%  "
$ws.Range("D3").Value = "1.674.10"
$ws.Range("E3").Value = "  +0.14%  "
$ws.Range("E4").Value = "  +0.26%  "
$ws.Range("D5").Value = "'214.95"
$ws.Range("E5").Value = "  -0.84%  "
$ws.Range("D6").Value = "'0.517"
$ws.Range("E6").Value = "  +1.46%  "
$ws.Range("E7").Value = "  +0.19%  "
$ws.Range("E8").Value = "  -0.42%  "
$ws.Range("E9").Value = "  -0.18%  "
$ws.Range("D10").Value = "'20.34"
$ws.Range("E10").Value = "  +0.61%  "
$ws.Range("E11").Value = "  +0.23%  "
$ws.Range("D12").Value = "1.908.83"
$ws.Range("E12").Value = "  +0.65%  "
$ws.Range("D13").Value = "1.642.14"
$ws.Range("E13").Value = "  -1.92%  "
$ws.Range("E14").Value = "  -0.49%  "
$ws.Range("E15").Value = "  -0.01%  "
$ws.Range("D16").Value = "'65.63"
$ws.Range("E16").Value = "  -0.40%  "
$ws.Range("D17").Value = "26.955.06"
$ws.Range("E17").Value = "  -0.54%  "
$ws.Range("D18").Value = "'8.11"
$ws.Range("E18").Value = "  +4.17%  "
$ws.Range("D19").Value = "'235.59"
$ws.Range("E19").Value = "  -0.92%  "
$ws.Range("E20").Value = "  -0.79%  "
$ws.Range("E21").Value = "  +0.14%  "
$ws.Range("E22").Value = "  -0.36%  "
$ws.Range("D23").Value = "'9.20"
$ws.Range("E23").Value = "  -1.48%  "
$ws.Range("E24").Value = "  -2.06%  "
$ws.Range("D25").Value = "'145.38"
$ws.Range("E25").Value = "  -0.18%  "
$ws.Range("D26").Value = "'7.19"
$ws.Range("E26").Value = "  +0.25%  "
$ws.Range("D27").Value = "'16.01"
$ws.Range("E27").Value = "  +0.31%  "
$ws.Range("E28").Value = "  -1.41%  "
$ws.Range("E29").Value = "  +0.31%  "
$ws.Range("E30").Value = "  -0.50%  "
$ws.Range("E31").Value = "  -0.45%  "
$ws.Range("E32").Value = "  +0.21%  "
$ws.Range("D33").Value = "1.477.01"
$ws.Range("E33").Value = "  -5.28%  "
$ws.Range("E34").Value = "  +0.71%  "
$ws.Range("E35").Value = "  +2.46%  "
$ws.Range("E36").Value = "  -0.27%  "
$ws.Range("E37").Value = "  +1.17%  "
$ws.Range("D38").Value = "'0.896"
$ws.Range("E38").Value = "  -1.40%  "
$ws.Range("E39").Value = "  +0.07%  "
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").Value = "'5.85"
$ws.Range("E40").Value = "  -4.00%  "
$ws.Range("B41").Value = "WEMIXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D41").Value = "'1.04"
$ws.Range("E41").Value = "  +6.41%  "
$ws.Range("E42").Value = "  +0.16%  "
$ws.Range("E43").Value = "  -0.16%  "
$ws.Range("D44").Value = "'67.34"
$ws.Range("E44").Value = "  +0.24%  "
$ws.Range("D45").Value = "1.814.79"
$ws.Range("E45").Value = "  +0.25%  "
$ws.Range("D46").Value = "'0.775"
$ws.Range("E46").Value = "  -0.80%  "
$ws.Range("D47").Value = "'90.64"
$ws.Range("E47").Value = "  -0.07%  "
$ws.Range("E48").Value = "  -0.71%  "
$ws.Range("E49").Value = "  +1.02%  "
$ws.Range("E50").Value = "  +0.13%  "
$ws.Range("E51").Value = "  -0.87%  "